$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new data row to Sheet1
$ws1.Range("A2").Value = "Western University of Health Sciences"
$ws1.Range("B2").Value = "Harriet K. & Philip Pumerantz Library"
$ws1.Range("C2").Value = "http://www.westernu.edu/library/"

# Autofit columns to match the bestFit column widths seen in the target file
$ws1.Columns("A:C").AutoFit() | Out-Null

# Add two more (empty) worksheets at the end, in order: Sheet2, then Sheet3
$wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null

# Re-select A3 on Sheet1 and re-activate Sheet1 as the visible tab
$ws1.Activate()
$ws1.Range("A3").Select() | Out-Null
